$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell with default (unstyled) formatting, used to reset style
# after forcing text format on numeric-looking values (avoids Excel auto-
# converting strings like "608.97" into numeric values).
$defaultStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = "64.310.47"
$ws.Range("E2").Value = "  -3.60%  "
$ws.Range("D3").Value = "3.160.35"
$ws.Range("E3").Value = "  -2.94%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "608.97"
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.33"
$ws.Range("D6").Style = $defaultStyle
$ws.Range("E6").Value = "  -6.88%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "3.158.32"
$ws.Range("E8").Value = "  -3.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.527"
$ws.Range("D9").Style = $defaultStyle
$ws.Range("E9").Value = "  -3.95%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.152"
$ws.Range("D10").Style = $defaultStyle
$ws.Range("E10").Value = "  -6.24%  "
$ws.Range("E11").Value = "  -7.65%  "
$ws.Range("E12").Value = "  -5.98%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000251"
$ws.Range("D13").Style = $defaultStyle
$ws.Range("E13").Value = "  -7.59%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.61"
$ws.Range("D14").Style = $defaultStyle
$ws.Range("E14").Value = "  -9.39%  "
$ws.Range("D15").Value = "3.679.68"
$ws.Range("E15").Value = "  -2.95%  "
$ws.Range("D16").Value = "64.324.46"
$ws.Range("E16").Value = "  -3.65%  "
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("D18").Value = "3.161.01"
$ws.Range("E18").Value = "  -3.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.94"
$ws.Range("D19").Style = $defaultStyle
$ws.Range("E19").Value = "  -6.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "481.94"
$ws.Range("D20").Style = $defaultStyle
$ws.Range("E20").Value = "  -5.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.72"
$ws.Range("D21").Style = $defaultStyle
$ws.Range("E21").Value = "  -4.56%  "
$ws.Range("E22").Value = "  -5.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.77"
$ws.Range("D23").Style = $defaultStyle
$ws.Range("E23").Value = "  -4.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.69"
$ws.Range("D24").Style = $defaultStyle
$ws.Range("E24").Value = "  -7.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.78"
$ws.Range("D25").Style = $defaultStyle
$ws.Range("E25").Value = "  -3.35%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  -5.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.47"
$ws.Range("D28").Style = $defaultStyle
$ws.Range("E28").Value = "  -6.72%  "
$ws.Range("E29").Value = "  -8.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.74"
$ws.Range("D30").Style = $defaultStyle
$ws.Range("E30").Value = "  -1.60%  "
$ws.Range("E31").Value = "  -18.88%  "
$ws.Range("E32").Value = "  -5.44%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.29"
$ws.Range("D34").Style = $defaultStyle
$ws.Range("E34").Value = "  -6.54%  "
$ws.Range("E35").Value = "  -4.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.55"
$ws.Range("D36").Style = $defaultStyle
$ws.Range("E36").Value = "  -2.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.99"
$ws.Range("D37").Style = $defaultStyle
$ws.Range("E37").Value = "  -7.12%  "
$ws.Range("D38").Value = "0.0₃0727"
$ws.Range("E38").Value = "  -8.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "452.59"
$ws.Range("D39").Style = $defaultStyle
$ws.Range("E39").Value = "  -8.65%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.95"
$ws.Range("D40").Style = $defaultStyle
$ws.Range("E40").Value = "  -12.33%  "
$ws.Range("E41").Value = "  -7.50%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.44"
$ws.Range("D42").Style = $defaultStyle
$ws.Range("E42").Value = "  -4.57%  "
$ws.Range("E43").Value = "  -8.11%  "
$ws.Range("D44").Value = "2.843.43"
$ws.Range("E44").Value = "  -4.21%  "
$ws.Range("E45").Value = "  -9.26%  "
$ws.Range("E46").Value = "  -8.80%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.51"
$ws.Range("D47").Style = $defaultStyle
$ws.Range("E47").Value = "  -7.81%  "
$ws.Range("E49").Value = "  -6.93%  "
$ws.Range("E50").Value = "  -5.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "119.69"
$ws.Range("D51").Style = $defaultStyle
$ws.Range("E51").Value = "  -1.57%  "
